# Refresh the cryptocurrency listing (prices / 1h volume change, and two
# pairs of rows whose rank order swapped) with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "41.609.15";
    "E2" = "  +0.17%  ";
    "D3" = "2.460.24";
    "E3" = "  -1.27%  ";
    "D4" = "0.999";
    "E4" = "  +0.34%  ";
    "D5" = "314.70";
    "E5" = "  +0.45%  ";
    "D6" = "92.13";
    "E6" = "  -2.14%  ";
    "D7" = "0.547";
    "E7" = "  +0.21%  ";
    "E8" = "  +0.22%  ";
    "D9" = "0.511";
    "E9" = "  +2.75%  ";
    "D10" = "32.40";
    "E10" = "  -1.18%  ";
    "D11" = "0.0818";
    "E11" = "  +4.36%  ";
    "E12" = "  +0.46%  ";
    "D13" = "2.838.49";
    "E13" = "  -1.29%  ";
    "D14" = "6.85";
    "E14" = "  +0.15%  ";
    "D15" = "15.78";
    "E15" = "  +1.82%  ";
    "D16" = "2.454.02";
    "E16" = "  -0.73%  ";
    "D17" = "0.777";
    "E17" = "  +2.41%  ";
    "D18" = "41.617.78";
    "E18" = "  +0.02%  ";
    "D19" = "6.47";
    "E19" = "  +2.42%  ";
    "D20" = "0.0₃0943";
    "E20" = "  +2.39%  ";
    "D21" = "70.67";
    "E21" = "  -0.01%  ";
    "D22" = "11.32";
    "E22" = "  +1.44%  ";
    "D23" = "238.36";
    "E23" = "  +1.11%  ";
    "D24" = "2.72";
    "E24" = "  +0.21%  ";
    "B25" = "ImmutableX";
    "C25" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";
    "D25" = "1.91";
    "E25" = "  +0.44%  ";
    "B26" = "Dai";
    "C26" = "https://coinranking.com/coin/MoTuySvg7+dai-dai";
    "D26" = "1.00";
    "E26" = "  +0.01%  ";
    "D27" = "24.34";
    "E27" = "  -0.84%  ";
    "E28" = "  +0.96%  ";
    "D29" = "9.70";
    "E29" = "  +0.57%  ";
    "D30" = "35.12";
    "E30" = "  -3.32%  ";
    "D31" = "155.82";
    "E31" = "  +0.88%  ";
    "E32" = "  +0.78%  ";
    "E33" = "  -0.21%  ";
    "D34" = "0.0758";
    "E34" = "  +0.10%  ";
    "E35" = "  -0.93%  ";
    "D36" = "17.43";
    "E36" = "  -3.92%  ";
    "D37" = "2.89";
    "E37" = "  -3.42%  ";
    "E38" = "  +1.36%  ";
    "E39" = "  +1.58%  ";
    "E40" = "  -2.07%  ";
    "E41" = "  -4.21%  ";
    "E42" = "  +0.18%  ";
    "D43" = "1.975.99";
    "E43" = "  +1.32%  ";
    "B44" = "EnergySwap";
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";
    "D44" = "18.85";
    "E44" = "  -4.04%  ";
    "B45" = "VeChain";
    "C45" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";
    "D45" = "0.0282";
    "E45" = "  -0.61%  ";
    "D46" = "2.92";
    "E46" = "  -1.93%  ";
    "E47" = "  +2.34%  ";
    "D48" = "2.697.23";
    "E48" = "  -1.27%  ";
    "D49" = "96.76";
    "D50" = "66.80";
    "E50" = "  -0.63%  ";
    "D51" = "52.40";
    "E51" = "  +2.95%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force the assignment to be stored as text (matching the source data,
    # which keeps prices/percentages as strings) instead of letting Excel's
    # automatic type inference coerce number-looking strings into doubles,
    # then restore the cell to its original (default/"Normal") style so no
    # stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells in the cryptos list."
